# [Outlook] (sensitivity label) Remove sensitivity label snippets from preview
#
# The "Snippets" table lists one row per Office JS API member / snippet
# mapping. This change removes the rows that document the (removed)
# sensitivity-label APIs: AppointmentCompose.sensitivityLabel (get/set),
# MessageCompose.sensitivityLabel (get/set), SensitivityLabel.getAsync,
# SensitivityLabel.setAsync, SensitivityLabelsCatalog.getAsync and
# SensitivityLabelsCatalog.getIsEnabledAsync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Rows (1-based, including the header row) that correspond to the
# sensitivity-label snippet entries being removed. Deleted from the
# bottom up so earlier row numbers remain valid while iterating.
$rowsToDelete = @(283, 282, 281, 280, 202, 201, 47, 46)

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).EntireRow.Delete()
}

# Restore the view state (selection) that was active after the edit.
$ws.Range("E273").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 248
